# Scheduled Universalis market-data refresh for the crafting-class leve-profit sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the derived LeveProfit(NQ/HQ)
# columns (H:N) with freshly pulled market prices. Cells that have no corresponding
# profit figure for the new price set are cleared instead of left stale.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 856.3333
$ws.Range("I40").Value = 684.5
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 684.5
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -509.5
$ws.Range("N40").Value = -1550
$ws.Range("H116").Value = 6000
$ws.Range("I116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558
$ws.Range("H135").Value = 832
$ws.Range("I135").Value = 832
$ws.Range("K135").Value = 7488
$ws.Range("M135").Value = -4953

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2720.7
$ws.Range("I2").Value = 1401
$ws.Range("K2").Value = 1401
$ws.Range("M2").Value = -1288
$ws.Range("H61").Value = 5000
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5424
$ws.Range("H63").Value = 1213.5385
$ws.Range("I63").Value = 1152.4
$ws.Range("J63").Value = 1417.3334
$ws.Range("K63").Value = 1152.4
$ws.Range("L63").Value = 1417.3334
$ws.Range("M63").Value = -466.4000000000001
$ws.Range("N63").Value = -2789.3334
$ws.Range("H66").Value = 1213.5385
$ws.Range("I66").Value = 1152.4
$ws.Range("J66").Value = 1417.3334
$ws.Range("K66").Value = 5762
$ws.Range("L66").Value = 7086.666999999999
$ws.Range("M66").Value = -2330
$ws.Range("N66").Value = -13950.667
$ws.Range("H97").Value = 965.5
$ws.Range("I97").Value = 965.5
$ws.Range("K97").Value = 965.5
$ws.Range("M97").Value = -469.5
$ws.Range("H102").Value = 624.75
$ws.Range("I102").Value = 624.75
$ws.Range("K102").Value = 624.75
$ws.Range("M102").Value = 997.25
$ws.Range("H116").Value = 2720.7
$ws.Range("I116").Value = 1401
$ws.Range("K116").Value = 1401
$ws.Range("M116").Value = 893
$ws.Range("H136").Value = 5000
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2720.7
$ws.Range("I3").Value = 1401
$ws.Range("K3").Value = 1401
$ws.Range("M3").Value = -1287
$ws.Range("H12").Value = 250.4
$ws.Range("I12").Value = 135.42857
$ws.Range("K12").Value = 135.42857
$ws.Range("M12").Value = 32.57142999999999
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H86").Value = 2592.0715
$ws.Range("I86").Value = 2158.8
$ws.Range("J86").Value = 2832.7778
$ws.Range("K86").Value = 2158.8
$ws.Range("L86").Value = 2832.7778
$ws.Range("M86").Value = -1035.8
$ws.Range("N86").Value = -5078.7778
$ws.Range("H89").Value = 2592.0715
$ws.Range("I89").Value = 2158.8
$ws.Range("J89").Value = 2832.7778
$ws.Range("K89").Value = 10794
$ws.Range("L89").Value = 14163.889
$ws.Range("M89").Value = -5178
$ws.Range("N89").Value = -25395.889
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 1955.4
$ws.Range("I107").Value = 1150
$ws.Range("K107").Value = 1150
$ws.Range("M107").Value = 770
$ws.Range("H122").Value = 1979899
$ws.Range("I122").Value = 1979899
$ws.Range("K122").Value = 1979899
$ws.Range("M122").Value = -1974999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 96.75
$ws.Range("I7").Value = 96.75
$ws.Range("K7").Value = 96.75
$ws.Range("M7").Value = 16.25
$ws.Range("H31").Value = 3399.6667
$ws.Range("I31").Value = 2442.0833
$ws.Range("J31").Value = 7230
$ws.Range("K31").Value = 2442.0833
$ws.Range("L31").Value = 7230
$ws.Range("M31").Value = -2147.0833
$ws.Range("N31").Value = -7820
$ws.Range("H34").Value = 3399.6667
$ws.Range("I34").Value = 2442.0833
$ws.Range("J34").Value = 7230
$ws.Range("K34").Value = 2442.0833
$ws.Range("L34").Value = 7230
$ws.Range("M34").Value = -2240.0833
$ws.Range("N34").Value = -7634
$ws.Range("H58").Value = 2531.3333
$ws.Range("I58").Value = 2437.6
$ws.Range("K58").Value = 2437.6
$ws.Range("M58").Value = -2234.6
$ws.Range("H107").Value = 1914.6666
$ws.Range("I107").Value = 1914.6666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1914.6666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 5.333399999999983
$ws.Range("N107").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H134").Value = 2533.2
$ws.Range("I134").Value = 2313.889
$ws.Range("K134").Value = 6941.667
$ws.Range("M134").Value = -4406.667
$ws.Range("H136").Value = 2531.3333
$ws.Range("I136").Value = 2437.6
$ws.Range("K136").Value = 7312.799999999999
$ws.Range("M136").Value = -4762.799999999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 493.4
$ws.Range("I5").Value = 493.4
$ws.Range("K5").Value = 1480.2
$ws.Range("M5").Value = -1368.2
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H34").Value = 714.7692
$ws.Range("J34").Value = 839.2
$ws.Range("L34").Value = 2517.6
$ws.Range("N34").Value = -2685.6
$ws.Range("H39").Value = 3700
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588
$ws.Range("H55").Value = 883.3333
$ws.Range("J55").Value = 925
$ws.Range("L55").Value = 2775
$ws.Range("N55").Value = -3129
$ws.Range("H135").Value = 493.4
$ws.Range("I135").Value = 493.4
$ws.Range("K135").Value = 4440.599999999999
$ws.Range("M135").Value = -1905.599999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29456
$ws.Range("J15").Value = 29456
$ws.Range("L15").Value = 29456
$ws.Range("N15").Value = -30032
$ws.Range("H23").Value = 900
$ws.Range("J23").Value = 900
$ws.Range("L23").Value = 900
$ws.Range("N23").Value = -1346
$ws.Range("H81").Value = 29456
$ws.Range("J81").Value = 29456
$ws.Range("L81").Value = 29456
$ws.Range("N81").Value = -31452
$ws.Range("H84").Value = 29456
$ws.Range("J84").Value = 29456
$ws.Range("L84").Value = 88368
$ws.Range("N84").Value = -98352
$ws.Range("H92").Value = 5927.857
$ws.Range("J92").Value = 5927.857
$ws.Range("L92").Value = 5927.857
$ws.Range("N92").Value = -9671.857
$ws.Range("H122").Value = 10419744
$ws.Range("J122").Value = 7335.6665
$ws.Range("L122").Value = 22006.9995
$ws.Range("N122").Value = -26906.9995

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -4887
$ws.Range("N4").ClearContents()
$ws.Range("H13").Value = 13002.333
$ws.Range("I13").Value = 9000
$ws.Range("J13").Value = 15003.5
$ws.Range("K13").Value = 9000
$ws.Range("L13").Value = 15003.5
$ws.Range("M13").Value = -8860
$ws.Range("N13").Value = -15283.5
$ws.Range("H16").Value = 683.3333
$ws.Range("I16").Value = 683.3333
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 683.3333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -513.3333
$ws.Range("N16").ClearContents()
$ws.Range("H20").Value = 5833.3335
$ws.Range("I20").Value = 3666.6667
$ws.Range("J20").Value = 8000
$ws.Range("K20").Value = 3666.6667
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = -3440.6667
$ws.Range("N20").Value = -8452
$ws.Range("H28").Value = 5000
$ws.Range("I28").Value = 5000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -4768
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 5000
$ws.Range("I37").Value = 5000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 5000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -4893
$ws.Range("N37").ClearContents()
$ws.Range("H46").Value = 3578.0303
$ws.Range("I46").Value = 2805
$ws.Range("K46").Value = 2805
$ws.Range("M46").Value = -2617
$ws.Range("H93").Value = 1500
$ws.Range("J93").Value = 1500
$ws.Range("L93").Value = 1500
$ws.Range("N93").Value = -3996

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 3005.5
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 5011
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 5011
$ws.Range("M20").Value = -760
$ws.Range("N20").Value = -5491
$ws.Range("H62").Value = 4654.7
$ws.Range("I62").Value = 4685.2856
$ws.Range("K62").Value = 4685.2856
$ws.Range("M62").Value = -4061.2856
$ws.Range("H65").Value = 4654.7
$ws.Range("I65").Value = 4685.2856
$ws.Range("K65").Value = 23426.428
$ws.Range("M65").Value = -20306.428
$ws.Range("H107").Value = 313.8
$ws.Range("I107").Value = 267.25
$ws.Range("K107").Value = 801.75
$ws.Range("M107").Value = 1118.25
$ws.Range("H132").Value = 1010.75
$ws.Range("I132").Value = 1093.5454
$ws.Range("J132").Value = 100
$ws.Range("K132").Value = 3280.6362
$ws.Range("L132").Value = 300
$ws.Range("M132").Value = -750.6361999999999
$ws.Range("N132").Value = -5360
$ws.Range("H136").Value = 2714.2856
$ws.Range("I136").Value = 2714.2856
$ws.Range("K136").Value = 8142.8568
$ws.Range("M136").Value = -5592.8568

